$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9293811321258545
$ws.Range("B1").Value = 1.368552088737488
$ws.Range("C1").Value = 1.780435085296631
$ws.Range("D1").Value = 4.868771076202393
$ws.Range("E1").Value = 4.589833736419678
